# Updates the "optimization_parameters" sheet to match the current format
# used in beta: split the old "Model" row into a "production_function"
# label row plus a new "L_curve" parameter row, and drop the obsolete
# "Deletion" row. Also moves the active-sheet/selection to reflect the
# new layout (optimization_parameters tab, row 17 selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

# Insert a new row right after the production_function/Sigmoid row (row 8)
# for the new L_curve parameter.
$ws.Rows.Item(9).Insert()

$ws.Cells.Item(8, 1).Value = "production_function"
$ws.Cells.Item(9, 1).Value = "L_curve"
$ws.Cells.Item(9, 2).Value = 0
$ws.Cells.Item(9, 2).NumberFormat = "0.00E+00"

# The old "Deletion" row (now shifted down to row 17) is no longer used;
# remove it so everything below shifts back up.
$ws.Rows.Item(17).Delete()

# Reflect the new active sheet / selection state.
$ws.Activate()
$ws.Rows.Item(17).Select() | Out-Null
